$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "report generated" timestamp in the title cell
$ws.Range("A1").Value = "Reporte generado: 2019-10-15 15:35:31"

# Build out two more detail rows (6 and 7) by copying the formatting
# (borders + fills) of the existing data row 5 down into them.
$ws.Range("A5:E5").Copy($ws.Range("A6:E6"))
$ws.Range("A5:E5").Copy($ws.Range("A7:E7"))

# Row 5: PEREIRA / SEDEPEREIRA1 / ELISA / 2019-10-15 13:15:54 / DESCONECTADO
$ws.Range("A5").Value = "PEREIRA"
$ws.Range("B5").Value = "SEDEPEREIRA1"
$ws.Range("C5").Value = "ELISA"
$ws.Range("D5").Value = "2019-10-15 13:15:54"
$ws.Range("E5").Value = "DESCONECTADO"

# Row 6: PEREIRA / SEDEPEREIRA1 / CAMILACAMARA / 2019-10-15 13:22:07 / DESCONECTADO
$ws.Range("A6").Value = "PEREIRA"
$ws.Range("B6").Value = "SEDEPEREIRA1"
$ws.Range("C6").Value = "CAMILACAMARA"
$ws.Range("D6").Value = "2019-10-15 13:22:07"
$ws.Range("E6").Value = "DESCONECTADO"

# Row 7: PEREIRA / SEDEPEREIRA1 / JESUS CAMARA / 2019-10-15 13:22:32 / DESCONECTADO
$ws.Range("A7").Value = "PEREIRA"
$ws.Range("B7").Value = "SEDEPEREIRA1"
$ws.Range("C7").Value = "JESUS CAMARA"
$ws.Range("D7").Value = "2019-10-15 13:22:32"
$ws.Range("E7").Value = "DESCONECTADO"

# Recolor the ESTADO column (now "DESCONECTADO") to red for every data row
$ws.Range("E5:E7").Interior.Color = 3618773

# Zebra-stripe the middle data row (row 6) with a light grey fill
$ws.Range("A6:D6").Interior.Color = 15921906
